$d = $word.ActiveDocument

$pairs = @(
    @("36×86=", "71×53="),
    @("33×49=", "97×20="),
    @("91×30=", "21×42="),
    @("62×96=", "52×85="),
    @("66×45=", "37×64="),
    @("64×92=", "51×81="),
    @("12×88=", "19×91="),
    @("37×22=", "33×51="),
    @("14×55=", "29×94="),
    @("50×33=", "37×55="),
    @("72×56=", "57×23="),
    @("95×12=", "13×18="),
    @("96×39=", "38×37="),
    @("77×31=", "46×80="),
    @("26×99=", "75×51="),
    @("59×62=", "65×61="),
    @("84×71=", "24×79="),
    @("76×77=", "72×34="),
    @("77×16=", "61×44="),
    @("75×24=", "58×50="),
    @("95×15=", "26×11="),
    @("56×37=", "71×84="),
    @("81×25=", "64×61="),
    @("31×66=", "88×33="),
    @("47×49=", "63×59="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
